# Update cryptocurrency price/volume data (sheet1) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.364.64"
$ws.Range("E2").Value = "  -4.21%  "
$ws.Range("D3").Value = "2.615.63"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "3.071.86"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").Value = "58.303.48"
$ws.Range("E14").Value = "  -4.33%  "
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("D17").Value = "2.625.37"
$ws.Range("E17").Value = "  -7.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "336.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.92%  "
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.413"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("D28").Value = "0.0₃0786"
$ws.Range("E28").Value = "  -5.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("E34").Value = "  -4.64%  "
$ws.Range("E35").Value = "  -5.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.883"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.849"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("E39").Value = "  -7.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0965"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "268.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.73%  "
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.37%  "
$ws.Range("E47").Value = "  -3.18%  "
$ws.Range("D48").Value = "2.023.76"
$ws.Range("E48").Value = "  -5.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0228"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.55%  "
